$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("C10").Value = "Anderson J.-Des. Tec. M"
$ws.Range("D10").Value = "-"
$ws.Range("E10").Value = "João Bosco-Gestão Integr"

# Row 11
$ws.Range("C11").Value = "Anderson J.-Des. Tec. M"
$ws.Range("D11").Value = "-"
$ws.Range("E11").Value = "João Bosco-Gestão Integr"

# Row 12
$ws.Range("C12").Value = "Anderson J.-Des. Tec. M"
$ws.Range("D12").Value = "-"

# Row 14
$ws.Range("C14").Value = "Anderson J.-Des. Tec. M"
$ws.Range("D14").Value = "-"

# Row 15
$ws.Range("B15").Value = "-"
$ws.Range("C15").Value = "Anderson J.-Des. Tec. M"
$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "[-, Aline S. M.-T. M. Metalicos]"

# Row 16
$ws.Range("B16").Value = "-"
$ws.Range("C16").Value = "Anderson J.-Des. Tec. M"
$ws.Range("D16").Value = "-"
$ws.Range("E16").Value = "[-, Aline S. M.-T. M. Metalicos]"
